# Update "Last Accessed" dates (column E, rows 2-15) from 9/9/2020 (44083)
# to 9/17/2020 (44091), and move the active cell selection to H12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 5).Value2 = 44091
}

$ws.Range("H12").Select()
